# edit.ps1 - applies the "Pendientes" document update described by the
# commit "Actualizacion de pendientes y utilerias."
#
# Strategy: for every paragraph that needs a structural change (merged
# runs, removed proofErr markers, removed/added bookmark, removed
# lastRenderedPageBreak, changed highlight, changed/added text) we
# rebuild that paragraph's XML explicitly and push it in with
# Range.InsertXML so the resulting run/element structure matches the
# target precisely, instead of relying on Find/Replace's automatic
# (and not entirely predictable) run-merging behaviour.

$d = $word.ActiveDocument

function New-PkgXml([string]$innerBodyXml) {
    return '<?xml version="1.0" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $innerBodyXml + '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'
}

function Set-ParagraphXml($paragraphIndex, [string]$newParagraphXml) {
    $p = $d.Paragraphs($paragraphIndex)
    $rng = $p.Range
    $rng.InsertXML((New-PkgXml $newParagraphXml))
}

# ---------------------------------------------------------------------
# 1) "Generar un nuevo programa ..." bullet: yellow -> green highlight
#    on every run, and drop the _GoBack bookmark (it moves further down
#    the document, see step 4).
# ---------------------------------------------------------------------
$p4 = '<w:p w:rsidR="00BB0726" w:rsidRPr="00234F83" w:rsidRDefault="00F947AE" w:rsidP="00B31802">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:jc w:val="both"/>' +
  '<w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>Generar un nuevo programa para que todas las actividades externas se autoricen primero por gerencia</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> únicamente, nadie </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>más</w:t></w:r>' +
  '<w:r><w:rPr><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>y posteriormente se verán reflejadas en los usuarios correspondientes</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t>.</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 4 $p4

# ---------------------------------------------------------------------
# 2) "Para las actividades que no proceden ..." bullet: yellow -> green
#    highlight and the three runs collapse into a single run.
# ---------------------------------------------------------------------
$p5 = '<w:p w:rsidR="00F947AE" w:rsidRPr="00234F83" w:rsidRDefault="00BB0726" w:rsidP="000C5B69">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:jc w:val="both"/>' +
  '<w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr>' +
  '<w:t>Para las actividades que no proceden, es decir, que rechaza gerencia, deben marcarse como no procede y no se deben mostrar como pendientes de autorizar.</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 5 $p5

# ---------------------------------------------------------------------
# 3) "Agregar check ..." bullet: drop the proofErr spell-check markers
#    around "check" and collapse the three runs into one.
# ---------------------------------------------------------------------
$p8 = '<w:p w:rsidR="007B3C6A" w:rsidRPr="00885EE0" w:rsidRDefault="007B3C6A" w:rsidP="00885EE0">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:jc w:val="both"/>' +
  '<w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="24"/><w:highlight w:val="green"/></w:rPr>' +
  '<w:t>Agregar check para confirmar el visto de las notificaciones y que se guarde la hora del visto.</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 8 $p8

# ---------------------------------------------------------------------
# 4) "Las notificaciones por correo ..." bullet gets a new leading
#    sentence, is re-split across four runs, and now hosts the
#    _GoBack bookmark that used to sit at the end of bullet 1.
# ---------------------------------------------------------------------
$p9 = '<w:p w:rsidR="007B3C6A" w:rsidRPr="008251B2" w:rsidRDefault="007B3C6A" w:rsidP="00885EE0">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="10"/></w:numPr><w:jc w:val="both"/>' +
  '<w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Generar nuevo programa que corra en segundo plano para l</w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">as notificaciones por correo se enviarán todas </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">desde </w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">el servidor. </w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 9 $p9

# ---------------------------------------------------------------------
# 5) "Corregir error al prender la pc cuando ..." bullet: drop the
#    proofErr markers around "sql" and merge the trailing three runs
#    into one (the leading "Corregir..." and "cuando" runs are left
#    untouched).
# ---------------------------------------------------------------------
$p12 = '<w:p w:rsidR="00A10BE3" w:rsidRDefault="00A10BE3" w:rsidP="00885EE0">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Corregir error al prender la pc </w:t></w:r>' +
  '<w:r w:rsidR="00FE56D2"><w:t>cuando</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> inicia las notificaciones automáticas, le gana a la instancia sql y truena.</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 12 $p12

# ---------------------------------------------------------------------
# 6) "Agregar el nombre del equipo servidor ..." bullet: drop the
#    proofErr markers around "sql" and merge those two runs into one
#    (the trailing "sys21alien03-pc creo." run is left untouched).
# ---------------------------------------------------------------------
$p13 = '<w:p w:rsidR="00D571AF" w:rsidRDefault="00D571AF" w:rsidP="00885EE0">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="3"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Agregar el nombre del equipo servidor para hacer la conexión sql, ahorita está fija por código a </w:t></w:r>' +
  '<w:r w:rsidR="00180F48"><w:t>sys21alien03-pc creo.</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 13 $p13

# ---------------------------------------------------------------------
# 7) "Enviar correo." bullet is replaced by a longer sentence split
#    across five runs.
# ---------------------------------------------------------------------
$p22 = '<w:p w:rsidR="0021311B" w:rsidRDefault="0021311B" w:rsidP="00885EE0">' +
  '<w:pPr><w:pStyle w:val="Prrafodelista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="4"/></w:numPr><w:jc w:val="both"/></w:pPr>' +
  '<w:r><w:t xml:space="preserve">Quitar </w:t></w:r>' +
  '<w:r><w:t>envío</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> de correo, para eso es el nuevo </w:t></w:r>' +
  '<w:r><w:t>programa de fondo</w:t></w:r>' +
  '<w:r><w:t>.</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 22 $p22

# ---------------------------------------------------------------------
# 8) "Catálogos:" heading: drop the stray lastRenderedPageBreak marker.
# ---------------------------------------------------------------------
$p26 = '<w:p w:rsidR="003412B7" w:rsidRPr="003412B7" w:rsidRDefault="00C94797" w:rsidP="00885EE0">' +
  '<w:pPr><w:jc w:val="both"/><w:rPr><w:b/></w:rPr></w:pPr>' +
  '<w:r w:rsidRPr="003412B7"><w:rPr><w:b/></w:rPr><w:t>Catálogos</w:t></w:r>' +
  '<w:r w:rsidR="003412B7" w:rsidRPr="003412B7"><w:rPr><w:b/></w:rPr><w:t>:</w:t></w:r>' +
  '</w:p>'
Set-ParagraphXml 26 $p26

Write-Host "Done applying Pendientes updates."
